# Fill in the missing "PageNumber" value for the "grid" action row (row 8,
# column B) and leave that cell selected, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 0
$ws.Range("B8").Select()
